# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "last updated" timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 06:30"

# Row 5: India (rank 9)
$ws.Range("B5").Value = 4930236
$ws.Range("C5").Value = 3322
$ws.Range("D5").Value = 3859399
$ws.Range("E5").Value = 990029

# Row 20: Pakistan (rank 24)
$ws.Range("B20").Value = 302424
$ws.Range("C20").Value = 404
$ws.Range("D20").Value = 290261
$ws.Range("E20").Value = 5774
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 6389

# Row 129: Tailandia (rank 133)
$ws.Range("B129").Value = 3480
$ws.Range("C129").Value = 5
$ws.Range("D129").Value = 3315
$ws.Range("E129").Value = 107

# Rows 135 & 136: Birmania overtakes Sri Lanka in ranking order.
# Row 135 keeps rank 139 but now shows Birmania's updated figures.
$ws.Range("A135").Value = "Birmania"
$ws.Range("B135").Value = 3299
$ws.Range("C135").Value = 104
$ws.Range("D135").Value = 790
$ws.Range("E135").Value = 2477
$ws.Range("H135").Value = 32

# Row 136 keeps rank 140 but now shows Sri Lanka's (unchanged) figures.
$ws.Range("A136").Value = "Sri Lanka"
$ws.Range("B136").Value = 3262
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 3005
$ws.Range("E136").Value = 244
$ws.Range("H136").Value = 13

# Row 187: Butan (rank 191)
$ws.Range("B187").Value = 246
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 173
$ws.Range("E187").Value = 73
